$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estonia Meistriliiga")

# Row 71
$ws.Range("F71").Value = "JK Tammeka Tartu"
$ws.Range("G71").Value = "FC Flora Tallinn"
$ws.Range("J71").Value = "A"
$ws.Range("B71").Value = 6139072
$ws.Range("H71").Value = 1
$ws.Range("I71").Value = 2
$ws.Range("K71").Value = 9
$ws.Range("L71").Value = 7
$ws.Range("M71").Value = 1.166
$ws.Range("N71").Value = 7
$ws.Range("O71").Value = 6
$ws.Range("P71").Value = 1.25
$ws.Range("Q71").Value = 1.75
$ws.Range("R71").Value = 1.9
$ws.Range("S71").Value = 1.9
$ws.Range("T71").Value = 3
$ws.Range("U71").Value = 1.95
$ws.Range("V71").Value = 1.85
$ws.Range("W71").Value = -1
$ws.Range("X71").Value = -1
$ws.Range("Y71").Value = 0.25
$ws.Range("Z71").Value = 0.8999999999999999
$ws.Range("AA71").Value = -1
$ws.Range("AB71").Value = 0
$ws.Range("AC71").Value = -0

# Row 72
$ws.Range("F72").Value = "Parnu JK Vaprus"
$ws.Range("G72").Value = "JK Trans Narva"
$ws.Range("J72").Value = "H"
$ws.Range("B72").Value = 6139071
$ws.Range("H72").Value = 3
$ws.Range("I72").Value = 2
$ws.Range("K72").Value = 2.4
$ws.Range("L72").Value = 3.2
$ws.Range("M72").Value = 2.6
$ws.Range("N72").Value = 3
$ws.Range("O72").Value = 3.25
$ws.Range("P72").Value = 2.2
$ws.Range("Q72").Value = 0.25
$ws.Range("R72").Value = 1.825
$ws.Range("S72").Value = 1.975
$ws.Range("T72").Value = 2.5
$ws.Range("U72").Value = 1.875
$ws.Range("V72").Value = 1.925
$ws.Range("W72").Value = 2
$ws.Range("X72").Value = -1
$ws.Range("Y72").Value = -1
$ws.Range("Z72").Value = 0.825
$ws.Range("AA72").Value = -1
$ws.Range("AB72").Value = 0.875
$ws.Range("AC72").Value = -1

# Row 104
$ws.Range("F104").Value = "FC Kuressaare"
$ws.Range("G104").Value = "Parnu JK Vaprus"
$ws.Range("J104").Value = "H"
$ws.Range("B104").Value = 6533597
$ws.Range("H104").Value = 1
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 2.5
$ws.Range("L104").Value = 3.4
$ws.Range("M104").Value = 2.5
$ws.Range("N104").Value = 2.15
$ws.Range("O104").Value = 3.6
$ws.Range("P104").Value = 2.875
$ws.Range("Q104").Value = -0.25
$ws.Range("R104").Value = 1.95
$ws.Range("S104").Value = 1.85
$ws.Range("T104").Value = 2.75
$ws.Range("U104").Value = 1.95
$ws.Range("V104").Value = 1.85
$ws.Range("W104").Value = 1.15
$ws.Range("X104").Value = -1
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = 0.95
$ws.Range("AA104").Value = -1
$ws.Range("AB104").Value = -1
$ws.Range("AC104").Value = 0.8500000000000001

# Row 106
$ws.Range("F106").Value = "JK Tallinna Kalev"
$ws.Range("G106").Value = "JK Trans Narva"
$ws.Range("J106").Value = "H"
$ws.Range("B106").Value = 6537869
$ws.Range("H106").Value = 5
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 1.6
$ws.Range("L106").Value = 4
$ws.Range("M106").Value = 4.5
$ws.Range("N106").Value = 1.65
$ws.Range("O106").Value = 4
$ws.Range("P106").Value = 4.333
$ws.Range("Q106").Value = -0.75
$ws.Range("R106").Value = 1.8
$ws.Range("S106").Value = 2
$ws.Range("T106").Value = 2.75
$ws.Range("U106").Value = 1.9
$ws.Range("V106").Value = 1.9
$ws.Range("W106").Value = 0.6499999999999999
$ws.Range("X106").Value = -1
$ws.Range("Y106").Value = -1
$ws.Range("Z106").Value = 0.8
$ws.Range("AA106").Value = -1
$ws.Range("AB106").Value = 0.8999999999999999
$ws.Range("AC106").Value = -1

# Row 107
$ws.Range("F107").Value = "FC Flora Tallinn"
$ws.Range("G107").Value = "JK Nomme Kalju"
$ws.Range("J107").Value = "D"
$ws.Range("B107").Value = 6537957
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 1.4
$ws.Range("L107").Value = 4
$ws.Range("M107").Value = 7.5
$ws.Range("N107").Value = 1.5
$ws.Range("O107").Value = 4.2
$ws.Range("P107").Value = 5
$ws.Range("Q107").Value = -1
$ws.Range("R107").Value = 1.85
$ws.Range("S107").Value = 1.95
$ws.Range("T107").Value = 2.75
$ws.Range("U107").Value = 1.85
$ws.Range("V107").Value = 1.95
$ws.Range("W107").Value = -1
$ws.Range("X107").Value = 3.2
$ws.Range("Y107").Value = -1
$ws.Range("Z107").Value = -1
$ws.Range("AA107").Value = 0.95
$ws.Range("AB107").Value = -1
$ws.Range("AC107").Value = 0.95

# Row 115
$ws.Range("F115").Value = "JK Nomme Kalju"
$ws.Range("G115").Value = "JK Trans Narva"
$ws.Range("J115").Value = "H"
$ws.Range("B115").Value = 7919323
$ws.Range("H115").Value = 3
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 1.285
$ws.Range("L115").Value = 5.5
$ws.Range("M115").Value = 6.5
$ws.Range("N115").Value = 1.571
$ws.Range("O115").Value = 4.75
$ws.Range("P115").Value = 4.2
$ws.Range("Q115").Value = -1
$ws.Range("R115").Value = 1.925
$ws.Range("S115").Value = 1.875
$ws.Range("T115").Value = 2.75
$ws.Range("U115").Value = 1.875
$ws.Range("V115").Value = 1.925
$ws.Range("W115").Value = 0.571
$ws.Range("X115").Value = -1
$ws.Range("Y115").Value = -1
$ws.Range("Z115").Value = 0.925
$ws.Range("AA115").Value = -1
$ws.Range("AB115").Value = 0.4375
$ws.Range("AC115").Value = -0.5

# Row 116
$ws.Range("F116").Value = "FC Kuressaare"
$ws.Range("G116").Value = "FC Levadia Tallinn"
$ws.Range("J116").Value = "A"
$ws.Range("B116").Value = 7919322
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 6
$ws.Range("K116").Value = 11
$ws.Range("L116").Value = 6
$ws.Range("M116").Value = 1.166
$ws.Range("N116").Value = 15
$ws.Range("O116").Value = 8.5
$ws.Range("P116").Value = 1.125
$ws.Range("Q116").Value = 2.5
$ws.Range("R116").Value = 1.825
$ws.Range("S116").Value = 1.975
$ws.Range("T116").Value = 3.25
$ws.Range("U116").Value = 1.9
$ws.Range("V116").Value = 1.9
$ws.Range("W116").Value = -1
$ws.Range("X116").Value = -1
$ws.Range("Y116").Value = 0.125
$ws.Range("Z116").Value = -1
$ws.Range("AA116").Value = 0.9750000000000001
$ws.Range("AB116").Value = 0.8999999999999999
$ws.Range("AC116").Value = -1
